# Trade #11 closed at 2026-02-16 21:54:08 - leadlag UP +0.000%
# Append the new trade row to both the "All Trades" log sheet (row 12)
# and the strategy-specific "leadlag" sheet (row 11).

$wb = $excel.ActiveWorkbook

# ---- "All Trades" sheet: new row 12 ----
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(12, 1).Value = 11
# Leading apostrophe forces literal text so Excel doesn't auto-convert
# these into date/time serial numbers.
$wsAll.Cells.Item(12, 2).Value = "'2026-02-16"
$wsAll.Cells.Item(12, 3).Value = "'21:54:08"
$wsAll.Cells.Item(12, 4).Value = "leadlag"
$wsAll.Cells.Item(12, 5).Value = "UP"
$wsAll.Cells.Item(12, 6).Value = 68342.2
# Column G (Exit Price) stays blank - trade is still OPEN.
$wsAll.Cells.Item(12, 8).Value = "OPEN"
$wsAll.Cells.Item(12, 9).Value = 0
$wsAll.Cells.Item(12, 10).Value = 0
$wsAll.Cells.Item(12, 11).Value = 100
$wsAll.Cells.Item(12, 12).Value = 0.6445
$wsAll.Cells.Item(12, 13).Value = "Coinbase leading with 0.064% move"
# Column N (Exit Reason) stays blank - trade is still OPEN.
$wsAll.Cells.Item(12, 15).Value = 0

# ---- "leadlag" sheet: new row 11 ----
$wsLead = $wb.Worksheets.Item("leadlag")

$wsLead.Cells.Item(11, 1).Value = 11
$wsLead.Cells.Item(11, 2).Value = "'2026-02-16"
$wsLead.Cells.Item(11, 3).Value = "'21:54:08"
$wsLead.Cells.Item(11, 4).Value = "leadlag"
$wsLead.Cells.Item(11, 5).Value = "UP"
$wsLead.Cells.Item(11, 6).Value = 68342.2
# Column G (Exit Price) stays blank - trade is still OPEN.
$wsLead.Cells.Item(11, 8).Value = "OPEN"
$wsLead.Cells.Item(11, 9).Value = 0
$wsLead.Cells.Item(11, 10).Value = 0
$wsLead.Cells.Item(11, 11).Value = 100
$wsLead.Cells.Item(11, 12).Value = 0.6445
$wsLead.Cells.Item(11, 13).Value = "Coinbase leading with 0.064% move"
# Column N (Exit Reason) stays blank - trade is still OPEN.
$wsLead.Cells.Item(11, 15).Value = 0
